# Apply roster changes to "LOS Galacticos" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Terry Rozier"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Miami Heat"

# Row 7
$ws.Range("A7").Value = "Paul George"
$ws.Range("B7").Value = "SG,SF,PF"
$ws.Range("C7").Value = "Philadelphia 76ers"

# Row 12
$ws.Range("A12").Value = "Marcus Smart"
$ws.Range("C12").Value = "Memphis Grizzlies"

# Row 13
$ws.Range("A13").Value = "Tyus Jones"
$ws.Range("B13").Value = "PG"
$ws.Range("C13").Value = "Phoenix Suns"

# Row 14
$ws.Range("A14").Value = "Collin Sexton"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Utah Jazz"

# Row 15
$ws.Range("A15").Value = "Jayson Tatum"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Boston Celtics"

# Row 16
$ws.Range("A16").Value = "Jrue Holiday"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Boston Celtics"

# Row 18
$ws.Range("A18").Value = "Jaden Ivey"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Detroit Pistons"
